$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values would otherwise be auto-detected as numbers by Excel;
# force them to remain plain text so formatting (e.g. trailing zeros) is preserved.

$ws.Range('D2').Value = '66.080.95'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '3.319.01'
$ws.Range('E3').Value = '  -0.62%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.23'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.66'
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.651'
$ws.Range('E7').Value = '  +2.68%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '3.315.52'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('E10').Value = '  -2.99%  '
$ws.Range('E11').Value = '  +2.28%  '
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '3.893.06'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D15').Value = '66.133.23'
$ws.Range('E15').Value = '  -0.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.18'
$ws.Range('E16').Value = '  -3.14%  '
$ws.Range('E17').Value = '  -1.33%  '
$ws.Range('D18').Value = '3.309.39'
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '425.17'
$ws.Range('E19').Value = '  -2.81%  '
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.14'
$ws.Range('E21').Value = '  -3.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.38'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '71.61'
$ws.Range('E23').Value = '  -2.57%  '
$ws.Range('E24').Value = '  -0.09%  '
$ws.Range('E25').Value = '  +0.24%  '
$ws.Range('D26').Value = '3.461.12'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.512'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('E28').Value = '  +4.95%  '
$ws.Range('E29').Value = '  -3.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.89'
$ws.Range('E30').Value = '  -1.59%  '
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '22.40'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('E34').Value = '  +0.06%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.56'
$ws.Range('E36').Value = '  -3.25%  '
$ws.Range('E37').Value = '  -4.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '160.62'
$ws.Range('E38').Value = '  -2.18%  '
$ws.Range('E39').Value = '  -3.09%  '
$ws.Range('D40').Value = '2.863.85'
$ws.Range('E40').Value = '  +0.93%  '
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '26.34'
$ws.Range('E42').Value = '  -3.55%  '
$ws.Range('E43').Value = '  -5.02%  '
$ws.Range('E44').Value = '  -2.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '39.84'
$ws.Range('E45').Value = '  -0.95%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0659'
$ws.Range('E46').Value = '  -1.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.91'
$ws.Range('E47').Value = '  -5.00%  '
$ws.Range('E48').Value = '  -2.42%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.14'
$ws.Range('E49').Value = '  -5.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '312.54'
$ws.Range('E50').Value = '  -2.79%  '
$ws.Range('E51').Value = '  -1.16%  '
